$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-CellText($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value2 = $val
    $c.ClearFormats()
}

$ws.Cells.Item(2, 4).Value2 = "68.565.94"
$ws.Cells.Item(2, 5).Value2 = "  +1.18%  "
$ws.Cells.Item(3, 4).Value2 = "3.247.40"
$ws.Cells.Item(3, 5).Value2 = "  -0.64%  "
Set-CellText 5 4 "584.52"
$ws.Cells.Item(5, 5).Value2 = "  +0.74%  "
Set-CellText 6 4 "181.05"
$ws.Cells.Item(6, 5).Value2 = "  -0.92%  "
Set-CellText 7 4 "0.999"
$ws.Cells.Item(7, 5).Value2 = "  -0.06%  "
Set-CellText 8 4 "0.594"
$ws.Cells.Item(8, 5).Value2 = "  -1.02%  "
Set-CellText 9 4 "0.137"
$ws.Cells.Item(9, 5).Value2 = "  +2.74%  "
Set-CellText 10 4 "6.66"
$ws.Cells.Item(10, 5).Value2 = "  -1.43%  "
Set-CellText 11 4 "0.422"
$ws.Cells.Item(11, 5).Value2 = "  +0.93%  "
$ws.Cells.Item(12, 4).Value2 = "3.811.56"
$ws.Cells.Item(12, 5).Value2 = "  -0.53%  "
Set-CellText 13 4 "0.138"
$ws.Cells.Item(13, 5).Value2 = "  +0.09%  "
Set-CellText 14 4 "28.28"
$ws.Cells.Item(14, 5).Value2 = "  -1.29%  "
$ws.Cells.Item(15, 4).Value2 = "68.498.46"
$ws.Cells.Item(15, 5).Value2 = "  +1.16%  "
$ws.Cells.Item(16, 5).Value2 = "  +1.76%  "
$ws.Cells.Item(17, 4).Value2 = "3.282.52"
$ws.Cells.Item(17, 5).Value2 = "  +0.40%  "
Set-CellText 18 4 "5.82"
$ws.Cells.Item(18, 5).Value2 = "  -0.52%  "
Set-CellText 19 4 "13.47"
$ws.Cells.Item(19, 5).Value2 = "  -1.10%  "
Set-CellText 20 4 "394.70"
$ws.Cells.Item(20, 5).Value2 = "  +5.26%  "
Set-CellText 21 4 "7.65"
$ws.Cells.Item(21, 5).Value2 = "  -0.07%  "
$ws.Cells.Item(22, 5).Value2 = "  -0.04%  "
$ws.Cells.Item(23, 5).Value2 = "  -0.03%  "
Set-CellText 24 4 "0.515"
$ws.Cells.Item(24, 5).Value2 = "  +0.44%  "
Set-CellText 25 4 "0.0000119"
$ws.Cells.Item(25, 5).Value2 = "  -1.22%  "
Set-CellText 26 4 "0.189"
$ws.Cells.Item(26, 5).Value2 = "  +4.46%  "
Set-CellText 27 4 "9.61"
$ws.Cells.Item(27, 5).Value2 = "  -0.63%  "
$ws.Cells.Item(28, 5).Value2 = "  +0.06%  "
$ws.Cells.Item(29, 5).Value2 = "  -0.12%  "
Set-CellText 30 4 "5.68"
$ws.Cells.Item(30, 5).Value2 = "  -0.67%  "
Set-CellText 31 4 "22.99"
$ws.Cells.Item(31, 5).Value2 = "  +0.94%  "
Set-CellText 32 4 "7.11"
$ws.Cells.Item(32, 5).Value2 = "  +2.47%  "
$ws.Cells.Item(33, 5).Value2 = "  +0.02%  "
Set-CellText 34 4 "1.28"
$ws.Cells.Item(34, 5).Value2 = "  -0.09%  "
Set-CellText 35 4 "164.17"
$ws.Cells.Item(35, 5).Value2 = "  +0.62%  "
Set-CellText 36 4 "1.48"
$ws.Cells.Item(36, 5).Value2 = "  -1.82%  "
Set-CellText 37 4 "1.93"
$ws.Cells.Item(37, 5).Value2 = "  +4.35%  "
Set-CellText 38 4 "0.819"
$ws.Cells.Item(38, 5).Value2 = "  -3.94%  "
Set-CellText 39 4 "4.58"
$ws.Cells.Item(39, 5).Value2 = "  -2.32%  "
Set-CellText 40 4 "26.24"
$ws.Cells.Item(40, 5).Value2 = "  -2.27%  "
Set-CellText 41 4 "6.56"
$ws.Cells.Item(41, 5).Value2 = "  -4.38%  "
Set-CellText 42 4 "41.41"
$ws.Cells.Item(42, 5).Value2 = "  +1.17%  "
$ws.Cells.Item(43, 5).Value2 = "  -6.06%  "
Set-CellText 44 4 "0.0687"
$ws.Cells.Item(44, 5).Value2 = "  +0.51%  "
Set-CellText 45 4 "343.74"
$ws.Cells.Item(45, 5).Value2 = "  -3.90%  "
$ws.Cells.Item(46, 4).Value2 = "2.603.17"
$ws.Cells.Item(46, 5).Value2 = "  -3.96%  "
Set-CellText 47 4 "24.60"
$ws.Cells.Item(47, 5).Value2 = "  -3.54%  "
Set-CellText 48 4 "0.0281"
$ws.Cells.Item(48, 5).Value2 = "  -0.27%  "
$ws.Cells.Item(49, 2).Value2 = "Arweave"
$ws.Cells.Item(49, 3).Value2 = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-CellText 49 4 "31.78"
$ws.Cells.Item(49, 5).Value2 = "  +1.18%  "
$ws.Cells.Item(50, 2).Value2 = "Cosmos"
$ws.Cells.Item(50, 3).Value2 = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-CellText 50 4 "6.30"
$ws.Cells.Item(50, 5).Value2 = "  +1.98%  "
Set-CellText 51 4 "0.102"
$ws.Cells.Item(51, 5).Value2 = "  -1.04%  "
